$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 901
$ws.Range("I32").Value = 800
$ws.Range("K32").Value = 800
$ws.Range("M32").Value = -474

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2729.6128
$ws.Range("I98").Value = 1317.5416
$ws.Range("J98").Value = 7571
$ws.Range("K98").Value = 1317.5416
$ws.Range("L98").Value = 7571
$ws.Range("M98").Value = 180.4584
$ws.Range("N98").Value = -10567

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7986.4287
$ws.Range("I116").Value = 2202
$ws.Range("J116").Value = 11546.077
$ws.Range("K116").Value = 2202
$ws.Range("L116").Value = 11546.077
$ws.Range("M116").Value = 1240
$ws.Range("N116").Value = -18430.077

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2729.6128
$ws.Range("I122").Value = 1317.5416
$ws.Range("J122").Value = 7571
$ws.Range("K122").Value = 3952.6248
$ws.Range("L122").Value = 22713
$ws.Range("M122").Value = -1502.6248
$ws.Range("N122").Value = -27613

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 128041.875
$ws.Range("I132").Value = 155691.69
$ws.Range("J132").Value = 8226
$ws.Range("K132").Value = 467075.07
$ws.Range("L132").Value = 24678
$ws.Range("M132").Value = -464545.07
$ws.Range("N132").Value = -29738

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1294.8163
$ws.Range("I74").Value = 991.0833
$ws.Range("J74").Value = 2135.923
$ws.Range("K74").Value = 991.0833
$ws.Range("L74").Value = 2135.923
$ws.Range("M74").Value = -117.0833
$ws.Range("N74").Value = -3883.923

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1294.8163
$ws.Range("I77").Value = 991.0833
$ws.Range("J77").Value = 2135.923
$ws.Range("K77").Value = 4955.4165
$ws.Range("L77").Value = 10679.615
$ws.Range("M77").Value = -587.4165000000003
$ws.Range("N77").Value = -19415.615

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2904.375
$ws.Range("I122").Value = 1771.909
$ws.Range("J122").Value = 5395.8
$ws.Range("K122").Value = 5315.727000000001
$ws.Range("L122").Value = 16187.4
$ws.Range("M122").Value = -2865.727000000001
$ws.Range("N122").Value = -21087.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10003549
$ws.Range("J99").Value = 5265.5557
$ws.Range("L99").Value = 5265.5557
$ws.Range("N99").Value = -8261.555700000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2353.8823
$ws.Range("I122").Value = 1447.3334
$ws.Range("J122").Value = 2848.3635
$ws.Range("K122").Value = 4342.0002
$ws.Range("L122").Value = 8545.0905
$ws.Range("M122").Value = -1892.0002
$ws.Range("N122").Value = -13445.0905

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 10003549
$ws.Range("J126").Value = 5265.5557
$ws.Range("L126").Value = 15796.6671
$ws.Range("N126").Value = -20736.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1928.5714
$ws.Range("I17").Value = 800
$ws.Range("J17").Value = 3433.3333
$ws.Range("K17").Value = 2400
$ws.Range("L17").Value = 10299.9999
$ws.Range("M17").Value = -2231
$ws.Range("N17").Value = -10637.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 6520.4443
$ws.Range("I34").Value = 150
$ws.Range("J34").Value = 6895.1763
$ws.Range("K34").Value = 450
$ws.Range("L34").Value = 20685.5289
$ws.Range("M34").Value = -366
$ws.Range("N34").Value = -20853.5289

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 126.57143
$ws.Range("I38").Value = 52.5
$ws.Range("J38").Value = 156.2
$ws.Range("K38").Value = 157.5
$ws.Range("L38").Value = 468.6
$ws.Range("M38").Value = 189.5
$ws.Range("N38").Value = -1162.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 10282.8
$ws.Range("J39").Value = 10282.8
$ws.Range("L39").Value = 30848.4
$ws.Range("N39").Value = -31436.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4869.478
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 4869.478
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 14608.434
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -14962.434

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 3101
$ws.Range("I59").Value = 1202
$ws.Range("J59").Value = 5000
$ws.Range("K59").Value = 3606
$ws.Range("L59").Value = 15000
$ws.Range("M59").Value = -3066
$ws.Range("N59").Value = -16080

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 547.2941
$ws.Range("I113").Value = 532.03705
$ws.Range("J113").Value = 564.4583
$ws.Range("K113").Value = 1596.11115
$ws.Range("L113").Value = 1693.3749
$ws.Range("M113").Value = 573.8888499999998
$ws.Range("N113").Value = -6033.3749

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7693260.5
$ws.Range("I131").Value = 166667470
$ws.Range("J131").Value = 960.4516
$ws.Range("K131").Value = 500002410
$ws.Range("L131").Value = 2881.3548
$ws.Range("M131").Value = -499997370
$ws.Range("N131").Value = -12961.3548

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2778.5557
$ws.Range("I122").Value = 1512.4706
$ws.Range("J122").Value = 4930.9
$ws.Range("K122").Value = 4537.4118
$ws.Range("L122").Value = 14792.7
$ws.Range("M122").Value = -2087.4118
$ws.Range("N122").Value = -19692.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4604.857
$ws.Range("I7").Value = 2591.5454
$ws.Range("J7").Value = 6819.5
$ws.Range("K7").Value = 2591.5454
$ws.Range("L7").Value = 6819.5
$ws.Range("M7").Value = -2479.5454
$ws.Range("N7").Value = -7043.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5607.4546
$ws.Range("I122").Value = 3575.5
$ws.Range("J122").Value = 6768.5713
$ws.Range("K122").Value = 10726.5
$ws.Range("L122").Value = 20305.7139
$ws.Range("M122").Value = -8276.5
$ws.Range("N122").Value = -25205.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4604.857
$ws.Range("I126").Value = 2591.5454
$ws.Range("J126").Value = 6819.5
$ws.Range("K126").Value = 7774.6362
$ws.Range("L126").Value = 20458.5
$ws.Range("M126").Value = -5304.6362
$ws.Range("N126").Value = -25398.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2956.6216
$ws.Range("I136").Value = 1572.862
$ws.Range("K136").Value = 4718.586
$ws.Range("M136").Value = -2168.586

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1762.65
$ws.Range("I136").Value = 729.5
$ws.Range("J136").Value = 3312.375
$ws.Range("K136").Value = 2188.5
$ws.Range("L136").Value = 9937.125
$ws.Range("M136").Value = 361.5
$ws.Range("N136").Value = -15037.125
